$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add units, insert new Power (W) column E, shift old Elapsed Time to F ---
# F1 is a brand-new cell: clone E1 formatting (bold/border/alignment) onto it first.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 1).Value = "Time"
$ws.Cells.Item(1, 2).Value = "Temp (°C)"
$ws.Cells.Item(1, 3).Value = "Voltage (V)"
$ws.Cells.Item(1, 4).Value = "Current (mA)"
$ws.Cells.Item(1, 5).Value = "Power (W)"
$ws.Cells.Item(1, 6).Value = "Elapsed Time (s)"

# --- Data rows 2-33: A=Time text, B=Temp, C=Voltage, D=Current, E=Power(new), F=Elapsed Time ---
$data = @(
    @("0d0h31m5", 23.9, 0, 0, 0, 0.5986261367797852),
    @("0d0h31m5", 24, 0, 0, 0, 2.105473279953003),
    @("0d0h31m5", 24, 0, 0, 0, 3.608244657516479),
    @("0d0h31m5", 24, 0, 0, 0, 5.115900754928589),
    @("0d0h31m5", 24, 0, 0, 0, 6.618702411651611),
    @("0d0h31m5", 23.9, 0, 0, 0, 8.126445770263672),
    @("0d0h31m5", 24.1, 0.42, 0, 0.00009799999999999998, 9.629774332046509),
    @("0d0h32m0", 23.9, 0, 0.04, 0, 11.13293600082397),
    @("0d0h32m2", 24.1, 0.56, 0, 0.0001742222222222223, 12.63974785804749),
    @("0d0h32m3", 23.9, 0, 0.01, 0, 14.14341926574707),
    @("0d0h32m5", 24, 0, 0.01, 0, 15.6506404876709),
    @("0d0h32m6", 24, 0, 0, 0, 17.15355730056763),
    @("0d0h32m8", 24.1, 0, 0, 0, 18.65697908401489),
    @("0d0h32m9", 24, 0, 0, 0, 20.15951108932495),
    @("0d0h32m1", 24.1, 1.51, 0, 0.001266722222222222, 21.66722917556763),
    @("0d0h32m1", 23.8, 0.71, 0, 0.0002800555555555555, 23.1697199344635),
    @("0d0h32m1", 24, 1.64, 0, 0.001494222222222222, 24.67335438728333),
    @("0d0h32m1", 24, 0, 0, 0, 26.18044018745422),
    @("0d0h32m1", 23.9, 0, 0, 0, 27.6843421459198),
    @("0d0h32m1", 23.9, 0, 0.02, 0, 29.19188261032104),
    @("0d0h32m2", 23.9, 0, 0.02, 0, 30.69759321212769),
    @("0d0h32m2", 24, 3.08, 0, 0.005270222222222222, 32.20168089866638),
    @("0d0h32m2", 23.9, 0, 0.02, 0, 33.70922541618347),
    @("0d0h32m2", 24.1, 0, 0, 0, 35.21572542190552),
    @("0d0h32m2", 23.9, 0.78, 0, 0.000338, 36.71937704086304),
    @("0d0h32m2", 23.9, 0, 0, 0, 38.23027110099792),
    @("0d0h32m2", 23.9, 0, 0, 0, 39.73340916633606),
    @("0d0h32m3", 23.9, 0, 0, 0, 41.23675870895386),
    @("0d0h32m3", 23.9, 0, 0, 0, 42.74416780471802),
    @("0d0h32m3", 23.9, 0, 0, 0, 44.24688673019409),
    @("0d0h32m3", 23.9, 0, 0, 0, 45.75023293495178),
    @("0d0h32m3", 24, 0, 0, 0, 47.25340867042542)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
